# Apply the "Apps That Are Needed" update:
#  1. Insert a new "Update 11" paragraph right before the
#     "Apps That Are Needed (Gap-Filling Ideas)" heading.
#  2. Re-point the four <w:lastRenderedPageBreak/> markers: each one moves
#     from the run it currently sits in to the start of the run that
#     immediately follows it (Word re-paginated and the rendered page
#     break now lands one run later each time).

$d = $word.ActiveDocument

function Replace-ParagraphXml($para, [string]$innerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) "Automatic prioritization ..." gains lastRenderedPageBreak; the
#    following horizontal-rule paragraph (pict 56ADC565) loses it.
# ---------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("Automatic prioritization based on deadlines, soldier availability, or operational urgency.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1 = $hit.Paragraphs(1)
$p2 = $p1.Next()

Replace-ParagraphXml $p1 '<w:p w14:paraId="6735B5F0" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="008E2B9A" w:rsidP="008E2B9A"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Automatic prioritization based on deadlines, soldier availability, or operational urgency.</w:t></w:r></w:p>'

Replace-ParagraphXml $p2 '<w:p w14:paraId="1291B493" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="00287175" w:rsidP="008E2B9A"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:pict w14:anchorId="56ADC565"><v:rect id="_x0000_i1026" style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t" fillcolor="#a0a0a0" stroked="f"/></w:pict></w:r></w:p>'

# ---------------------------------------------------------------------
# 2) "Integration with GPS ..." gains lastRenderedPageBreak; "AI-driven
#    load plans ..." loses it.
# ---------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("Integration with GPS for convoy route optimization.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1 = $hit.Paragraphs(1)

$hit2 = $d.Content
$hit2.Find.Execute("AI-driven load plans based on historical unit data and mission profiles.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $hit2.Paragraphs(1)

Replace-ParagraphXml $p1 '<w:p w14:paraId="1394EBE7" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="008E2B9A" w:rsidP="008E2B9A"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Integration with GPS for convoy route optimization.</w:t></w:r></w:p>'

Replace-ParagraphXml $p2 '<w:p w14:paraId="29E5CB89" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="008E2B9A" w:rsidP="008E2B9A"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>AI-driven load plans based on historical unit data and mission profiles.</w:t></w:r></w:p>'

# ---------------------------------------------------------------------
# 3) "Voice-activated assistant ..." gains lastRenderedPageBreak;
#    "Conversational AI assistant ..." loses it.
# ---------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("Voice-activated assistant for looking up tasks or procedures.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1 = $hit.Paragraphs(1)

$hit2 = $d.Content
$hit2.Find.Execute("Conversational AI assistant for quick answers to doctrinal or procedural questions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $hit2.Paragraphs(1)

Replace-ParagraphXml $p1 '<w:p w14:paraId="486DEF10" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="008E2B9A" w:rsidP="008E2B9A"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Voice-activated assistant for looking up tasks or procedures.</w:t></w:r></w:p>'

Replace-ParagraphXml $p2 '<w:p w14:paraId="0E07F8A2" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="008E2B9A" w:rsidP="008E2B9A"><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Conversational AI assistant for quick answers to doctrinal or procedural questions.</w:t></w:r></w:p>'

# ---------------------------------------------------------------------
# 4) The horizontal-rule paragraph (pict 75C1F7DA) gains
#    lastRenderedPageBreak; "If you're looking to build ..." loses it.
# ---------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("If you're looking to build one of these apps, the ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $hit.Paragraphs(1)
$p1 = $p2.Previous()

Replace-ParagraphXml $p1 '<w:p w14:paraId="735FB0D3" w14:textId="77777777" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="00287175" w:rsidP="008E2B9A"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:pict w14:anchorId="75C1F7DA"><v:rect id="_x0000_i1035" style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t" fillcolor="#a0a0a0" stroked="f"/></w:pict></w:r></w:p>'

Replace-ParagraphXml $p2 '<w:p w14:paraId="22C611E1" w14:textId="7F340D89" w:rsidR="008E2B9A" w:rsidRPr="008E2B9A" w:rsidRDefault="008E2B9A" w:rsidP="008E2B9A"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">If you''re looking to build one of these apps, the </w:t></w:r><w:r w:rsidRPr="008E2B9A"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Platoon Leader Dashboard</w:t></w:r><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> or the </w:t></w:r><w:r w:rsidRPr="008E2B9A"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>AI-Powered Task Management Tool</w:t></w:r><w:r w:rsidRPr="008E2B9A"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> would likely have the broadest impact. Let me know which idea you''d like to explore further!</w:t></w:r></w:p>'

# ---------------------------------------------------------------------
# 5) Insert the new "Update 11" paragraph right before the
#    "Apps That Are Needed (Gap-Filling Ideas)" heading (same run/
#    paragraph formatting as the other "Update N" lines above it).
# ---------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("Apps That Are Needed (Gap-Filling Ideas)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $hit.Paragraphs(1)
$headingPara.Range.InsertParagraphBefore()

$hit2 = $d.Content
$hit2.Find.Execute("Apps That Are Needed (Gap-Filling Ideas)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara2 = $hit2.Paragraphs(1)
$newPara = $headingPara2.Previous()
$newPara.Range.Text = "Update 11"
$newPara.Range.Font.Bold = -1
$newPara.Range.Font.Size = 11

Write-Output "done"
